# "Generate Report for Archive"
# - Update the Status value from "Ready for handoff" to "In Translation"
#   (this shared string is referenced by the Overview sheet's zh-cn/de-de
#   status columns as well as each language sheet's Status column).
# - Narrow the "Latest Handoff Datetime"-width columns (Overview!E:F and the
#   Status column on each language sheet) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: Ready for handoff -> In Translation ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Column width changes (closest reachable value to 13.4101845877511) ---
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
